# "fatti conti soglia dal basso (prima di ordinarli)"
#
# 1) Rename sheet 7 "soglia diodo" -> "soglia diodo crescenti"
# 2) Append raw (not-yet-sorted) threshold measurements in rows 22-40
#    on columns A (V) / B (I(uA)) of that sheet
# 3) Scroll the sheet view so row 16 is at the top (best effort - mirrors
#    the author having scrolled down to add these rows from below)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(7)

$ws.Name = "soglia diodo crescenti"

$newRows = @(
    @{ r = 22; a = 0.66300000000000003; b = 326.60000000000002 },
    @{ r = 23; a = 0.67400000000000004; b = 360.30000000000001 },
    @{ r = 24; a = 0.68400000000000005; b = 396.69999999999999 },
    @{ r = 25; a = 0.69399999999999995; b = 432.60000000000002 },
    @{ r = 26; a = 0.71399999999999997; b = 507.30000000000001 },
    @{ r = 27; a = 0.72499999999999998; b = 545.20000000000005 },
    @{ r = 28; a = 0.73499999999999999; b = 585.89999999999998 },
    @{ r = 29; a = 0.74099999999999999; b = 605 },
    @{ r = 30; a = 0.76600000000000001; b = 708 },
    @{ r = 31; a = 0.77600000000000002; b = 750 },
    @{ r = 32; a = 0.79100000000000004; b = 814 },
    @{ r = 33; a = 0.80700000000000005; b = 878 },
    @{ r = 34; a = 0.81699999999999995; b = 921 },
    @{ r = 35; a = 0.83199999999999996; b = 987 },
    @{ r = 36; a = 0.83699999999999997; b = 1009 },
    @{ r = 37; a = 0.84799999999999998; b = 1053 },
    @{ r = 38; a = 0.86799999999999999; b = 1143 },
    @{ r = 39; a = 0.88400000000000001; b = 47450 },
    @{ r = 40; a = 0.91400000000000003; b = 57120 }
)

foreach ($row in $newRows) {
    $ws.Cells.Item($row.r, 1).Value = $row.a
    $ws.Cells.Item($row.r, 2).Value = $row.b
    $ws.Rows.Item($row.r).RowHeight = 14.25
}

$ws.Activate()
$ws.Range("A1").Select()
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
